# Auto-generated: refresh market-price-derived columns (H:N) per scheduled runner update
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 3083.1667
$ws.Range("J127").Value = 4212.25
$ws.Range("L127").Value = 12636.75
$ws.Range("N127").Value = -22556.75
$ws.Range("H138").Value = 14495288
$ws.Range("I138").Value = 956.2143
$ws.Range("J138").Value = 24393856
$ws.Range("K138").Value = 2868.6429
$ws.Range("L138").Value = 73181568
$ws.Range("M138").Value = 2271.3571
$ws.Range("N138").Value = -73191848

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5009
$ws.Range("I2").Value = 4676.75
$ws.Range("K2").Value = 4676.75
$ws.Range("M2").Value = -4563.75
$ws.Range("H74").Value = 37740.574
$ws.Range("I74").Value = 47078.69
$ws.Range("J74").Value = 3056.1428
$ws.Range("K74").Value = 47078.69
$ws.Range("L74").Value = 3056.1428
$ws.Range("M74").Value = -46204.69
$ws.Range("N74").Value = -4804.1428
$ws.Range("H77").Value = 37740.574
$ws.Range("I77").Value = 47078.69
$ws.Range("J77").Value = 3056.1428
$ws.Range("K77").Value = 235393.45
$ws.Range("L77").Value = 15280.714
$ws.Range("M77").Value = -231025.45
$ws.Range("N77").Value = -24016.714
$ws.Range("H110").Value = 30086
$ws.Range("I110").Value = 32762
$ws.Range("K110").Value = 32762
$ws.Range("M110").Value = -30717
$ws.Range("H116").Value = 5009
$ws.Range("I116").Value = 4676.75
$ws.Range("K116").Value = 4676.75
$ws.Range("M116").Value = -2382.75
$ws.Range("H122").Value = 2228.4285
$ws.Range("I122").Value = 2015.2307
$ws.Range("K122").Value = 6045.6921
$ws.Range("M122").Value = -3595.6921

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5009
$ws.Range("I3").Value = 4676.75
$ws.Range("K3").Value = 4676.75
$ws.Range("M3").Value = -4562.75
$ws.Range("H22").Value = 350
$ws.Range("J22").Value = 450
$ws.Range("L22").Value = 450
$ws.Range("N22").Value = -796
$ws.Range("H86").Value = 13404.28
$ws.Range("I86").Value = 7047.0586
$ws.Range("J86").Value = 26913.375
$ws.Range("K86").Value = 7047.0586
$ws.Range("L86").Value = 26913.375
$ws.Range("M86").Value = -5924.0586
$ws.Range("N86").Value = -29159.375
$ws.Range("H89").Value = 13404.28
$ws.Range("I89").Value = 7047.0586
$ws.Range("J89").Value = 26913.375
$ws.Range("K89").Value = 35235.29300000001
$ws.Range("L89").Value = 134566.875
$ws.Range("M89").Value = -29619.29300000001
$ws.Range("N89").Value = -145798.875
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3166.4666
$ws.Range("I31").Value = 1454.5
$ws.Range("J31").Value = 5509.1577
$ws.Range("K31").Value = 1454.5
$ws.Range("L31").Value = 5509.1577
$ws.Range("M31").Value = -1159.5
$ws.Range("N31").Value = -6099.1577
$ws.Range("H34").Value = 3166.4666
$ws.Range("I34").Value = 1454.5
$ws.Range("J34").Value = 5509.1577
$ws.Range("K34").Value = 1454.5
$ws.Range("L34").Value = 5509.1577
$ws.Range("M34").Value = -1252.5
$ws.Range("N34").Value = -5913.1577
$ws.Range("H48").Value = 36000
$ws.Range("J48").Value = 36000
$ws.Range("L48").Value = 36000
$ws.Range("N48").Value = -36952
$ws.Range("H132").Value = 3294.04
$ws.Range("I132").Value = 3204.3171
$ws.Range("K132").Value = 9612.951300000001
$ws.Range("M132").Value = -7082.951300000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 51.5
$ws.Range("I14").Value = 51.5
$ws.Range("K14").Value = 154.5
$ws.Range("M14").Value = 18.5
$ws.Range("H94").Value = 10000.5
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 10000.5
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 30001.5
$ws.Range("N94").Value = -31353.5
$ws.Range("H99").Value = 3265.625
$ws.Range("J99").Value = 2500.5
$ws.Range("L99").Value = 7501.5
$ws.Range("N99").Value = -11993.5
$ws.Range("H125").Value = 7916.6665
$ws.Range("J125").Value = 10000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840
$ws.Range("H139").Value = 3405.6667
$ws.Range("I139").Value = 3226.4285
$ws.Range("K139").Value = 9679.2855
$ws.Range("M139").Value = -4539.2855
$ws.Range("M94").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 100.86667
$ws.Range("J2").Value = 261.75
$ws.Range("L2").Value = 261.75
$ws.Range("N2").Value = -487.75
$ws.Range("H52").Value = 34134.145
$ws.Range("I52").Value = 20000
$ws.Range("K52").Value = 20000
$ws.Range("M52").Value = -19741
$ws.Range("H97").Value = 895.4091
$ws.Range("I97").Value = 683.5
$ws.Range("J97").Value = 1149.7
$ws.Range("K97").Value = 683.5
$ws.Range("L97").Value = 1149.7
$ws.Range("M97").Value = -187.5
$ws.Range("N97").Value = -2141.7
$ws.Range("H126").Value = 4332.5
$ws.Range("I126").Value = 3666.6667
$ws.Range("K126").Value = 11000.0001
$ws.Range("M126").Value = -8530.000100000001
$ws.Range("H132").Value = 1240.8462
$ws.Range("I132").Value = 1240.8462
$ws.Range("K132").Value = 3722.5386
$ws.Range("M132").Value = -1192.5386

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5138.0835
$ws.Range("I7").Value = 3808.2856
$ws.Range("K7").Value = 3808.2856
$ws.Range("M7").Value = -3696.2856
$ws.Range("H16").Value = 843.9143
$ws.Range("I16").Value = 817.64514
$ws.Range("J16").Value = 1047.5
$ws.Range("K16").Value = 817.64514
$ws.Range("L16").Value = 1047.5
$ws.Range("M16").Value = -647.64514
$ws.Range("N16").Value = -1387.5
$ws.Range("H40").Value = 8335454
$ws.Range("I40").Value = 9805533
$ws.Range("J40").Value = 5004.3335
$ws.Range("K40").Value = 9805533
$ws.Range("L40").Value = 5004.3335
$ws.Range("M40").Value = -9805397
$ws.Range("N40").Value = -5276.3335
$ws.Range("H50").Value = 41746.25
$ws.Range("H82").Value = 2517.875
$ws.Range("I82").Value = 2352.25
$ws.Range("K82").Value = 2352.25
$ws.Range("M82").Value = -1991.25
$ws.Range("H85").Value = 2517.875
$ws.Range("I85").Value = 2352.25
$ws.Range("K85").Value = 2352.25
$ws.Range("M85").Value = -1104.25
$ws.Range("H87").Value = 41637
$ws.Range("I87").Value = 43849
$ws.Range("J87").Value = 35001
$ws.Range("K87").Value = 43849
$ws.Range("L87").Value = 35001
$ws.Range("M87").Value = -42726
$ws.Range("N87").Value = -37247
$ws.Range("H90").Value = 41637
$ws.Range("I90").Value = 43849
$ws.Range("J90").Value = 35001
$ws.Range("K90").Value = 131547
$ws.Range("L90").Value = 105003
$ws.Range("M90").Value = -125931
$ws.Range("N90").Value = -116235
$ws.Range("H122").Value = 5376.2856
$ws.Range("I122").Value = 4467.8
$ws.Range("K122").Value = 13403.4
$ws.Range("M122").Value = -10953.4
$ws.Range("H126").Value = 5138.0835
$ws.Range("I126").Value = 3808.2856
$ws.Range("K126").Value = 11424.8568
$ws.Range("M126").Value = -8954.856800000001
$ws.Range("H132").Value = 1987.7428
$ws.Range("I132").Value = 1799.4062
$ws.Range("K132").Value = 5398.2186
$ws.Range("M132").Value = -2868.2186
$ws.Range("H136").Value = 2686.6667
$ws.Range("I136").Value = 2612.9546
$ws.Range("J136").Value = 3497.5
$ws.Range("K136").Value = 7838.8638
$ws.Range("L136").Value = 10492.5
$ws.Range("M136").Value = -5288.8638
$ws.Range("N136").Value = -15592.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 35567.094
$ws.Range("J52").Value = 37772.723
$ws.Range("L52").Value = 37772.723
$ws.Range("N52").Value = -38224.723
$ws.Range("H81").Value = 1927
$ws.Range("I81").Value = 1391
$ws.Range("J81").Value = 2999
$ws.Range("K81").Value = 2782
$ws.Range("L81").Value = 5998
$ws.Range("M81").Value = -1721
$ws.Range("N81").Value = -8120
$ws.Range("H84").Value = 1927
$ws.Range("I84").Value = 1391
$ws.Range("J84").Value = 2999
$ws.Range("K84").Value = 13910
$ws.Range("L84").Value = 29990
$ws.Range("M84").Value = -8606
$ws.Range("N84").Value = -40598
$ws.Range("H136").Value = 2988.9106
$ws.Range("I136").Value = 2019.8649
$ws.Range("K136").Value = 6059.5947
$ws.Range("M136").Value = -3509.5947
$ws.Range("H139").Value = 120000
$ws.Range("J139").Value = 120000
$ws.Range("L139").Value = 120000
$ws.Range("N139").Value = -130280
